# v1.0.10 Add analyses for all three MCP1, MCP1_pg_ml_2015, and MCP1_pg_ug_2015.
# Add comparison between MCP1, MCP1_pg_ml_2015, and MCP1_pg_ug_2015.
# Add (and fixed) ordinal regression. Double checked which measurement to use.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Outcome / CpG labels used across the results table
$outcomeMajor    = "epmajor.90days"
$outcomeStroke   = "epstroke.90days"
$outcomeCoronary = "epcoronary.90days"
$outcomeCvdeath  = "epcvdeath.90days"

$cpgUg = "MCP1_pg_ug_2015_rank"
$cpgMl = "MCP1_pg_ml_2015_rank"
$cpgRk = "MCP1_rank"

$dataset = "AEDB.CEA"

# Full replacement data: Dataset, Outcome, CpG, Beta, s.e.m., HR, low95CI, up95CI, Z-value, P-value, SampleSize, N_events
$rows = @(
    @($dataset, $outcomeMajor,    $cpgUg, 1.108808975717,       0.381786002585533,  3.03074655122802,  1.43405831977462,   6.40519603081702,   2.90426827648976,   0.00368112584025718,  1029, 37),
    @($dataset, $outcomeMajor,    $cpgMl, 0.815177056352301,    0.373123633742323,  2.25957570916694,  1.08747091002143,   4.69500594306165,   2.18473712902158,   0.0289081261739204,   1029, 37),
    @($dataset, $outcomeMajor,    $cpgRk, -0.462104770413037,   0.536552389179315,  0.62995633569884,  0.220083677382947,  1.80315500725025,   -0.861248183275915, 0.389101367663012,    493,  17),
    @($dataset, $outcomeStroke,   $cpgUg, 1.10546122275982,     0.446587176506691,  3.02061732503239,  1.25878926323529,   7.24833718460176,   2.47535370676559,   0.0133104280099471,   1029, 27),
    @($dataset, $outcomeStroke,   $cpgMl, 0.924775855993483,    0.443385308912267,  2.52130306204764,  1.05732343189811,   6.01232218913255,   2.08571605194179,   0.0370043515312527,   1029, 27),
    @($dataset, $outcomeStroke,   $cpgRk, 0.164612344502335,    0.593397148988678,  1.17893600909127,  0.368451545376,     3.77224666574184,   0.2774066993461,    0.781467838081554,    493,  13),
    @($dataset, $outcomeCoronary, $cpgUg, 0.791264058554015,    0.623215020211065,  2.2061834094403,   0.650354694511652,  7.48398570374651,   1.2696485689418,    0.20420984052804,     1029, 13),
    @($dataset, $outcomeCoronary, $cpgMl, 1.29752697608952,     0.667509291082067,  3.66023362041083,  0.989266782061556,  13.5426665475079,   1.94383358168119,   0.0519155227223299,   1029, 13),
    @($dataset, $outcomeCoronary, $cpgRk, -1.36579275519938,    1.1180932199538,    0.255178301856539, 0.0285169832320325, 2.28341003704921,   -1.22153746291013,  0.221882590242073,    493,  5),
    @($dataset, $outcomeCvdeath,  $cpgUg, 2.91950274774125,     0.857554692296473,  18.532070054555,   3.45110873141841,   99.5151550515712,   3.4044507877661,    0.000662972450947386, 1029, 6),
    @($dataset, $outcomeCvdeath,  $cpgMl, 0.102224662788281,    0.927060709388879,  1.10763228753517,  0.179997138108792,  6.81593772701468,   0.110267495702269,  0.91219723470555,     1029, 6),
    @($dataset, $outcomeCvdeath,  $cpgRk, 31.2435242613778,     40.3755171219442,   37058700128178.8,  0.00000000000000000000158692416495058, 865414545649102000861002461397981454126439989248.0, 0.773823507127217, 0.439035169310639, 493, 3)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value  = $row[0]
    $ws.Cells.Item($r, 2).Value  = $row[1]
    $ws.Cells.Item($r, 3).Value  = $row[2]
    $ws.Cells.Item($r, 4).Value  = $row[3]
    $ws.Cells.Item($r, 5).Value  = $row[4]
    $ws.Cells.Item($r, 6).Value  = $row[5]
    $ws.Cells.Item($r, 7).Value  = $row[6]
    $ws.Cells.Item($r, 8).Value  = $row[7]
    $ws.Cells.Item($r, 9).Value  = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $r++
}
